# feat: add logviewer in alarm page
#
# - Rows 31/32 ("  HISTORIQUE DES ALARMES" / "  ALARMS HISTORY" and
#   "  ALARMES ACTIVES" / "  ACTIVE ALARMS", both indented sub-items) are
#   repurposed into top-level, un-indented "HISTORIQUE DES ALARMES" /
#   "ALARMES ACTIVES" labels (column A) with their English counterparts
#   "ALARMS HISTORY" / "ACTIVE ALARMS" (column B).
# - Two brand new rows are appended at the bottom of the translation table
#   for the new log viewer: "JOURNAL D'EVENEMENTS" / "EVENT LOG" and
#   "HISTORIQUES" / "HISTORY".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strings for the event-log feature are written first so that they end
# up ordered before the (reused) alarm strings in the shared-strings table,
# matching the order the original author's workbook ended up with.
$ws.Range("A215").Value = "JOURNAL D'EVENEMENTS"
$ws.Range("B215").Value = "EVENT LOG"

# Repurpose rows 31 & 32: drop the leading double-space indent and give
# them their own place as stand-alone labels for the alarm history /
# active alarm pages (French column first, then English column).
$ws.Range("A31").Value = "HISTORIQUE DES ALARMES"
$ws.Range("A32").Value = "ALARMES ACTIVES"
$ws.Range("B31").Value = "ALARMS HISTORY"
$ws.Range("B32").Value = "ACTIVE ALARMS"

# Append the final new pair of rows for the generic "History" section of
# the log viewer.
$ws.Range("A216").Value = "HISTORIQUES"
$ws.Range("B216").Value = "HISTORY"

# Match the author's final cursor position / selection on the sheet.
$ws.Range("B217").Select()
